$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the list of cell updates (address, new text value) taken from the commit diff.
# Values are written as text (NumberFormat "@") so Excel keeps them as literal strings
# instead of re-interpreting them as numbers/percentages, matching the original inlineStr cells.
$updates = @(
    @{ Cell = 'D2'; Value = '298.77' }
    @{ Cell = 'E2'; Value = '0.78%' }
    @{ Cell = 'D3'; Value = '31.20' }
    @{ Cell = 'E3'; Value = '-0.08%' }
    @{ Cell = 'D4'; Value = '5.154' }
    @{ Cell = 'E4'; Value = '0.93%' }
    @{ Cell = 'D5'; Value = '0.08032' }
    @{ Cell = 'E5'; Value = '9.47%' }
    @{ Cell = 'D6'; Value = '2.719' }
    @{ Cell = 'E6'; Value = '62.50%' }
    @{ Cell = 'D7'; Value = '7.824' }
    @{ Cell = 'D8'; Value = '3.825' }
    @{ Cell = 'E8'; Value = '2.30%' }
    @{ Cell = 'D9'; Value = '0.9156' }
    @{ Cell = 'E9'; Value = '-0.18%' }
    @{ Cell = 'D10'; Value = '0.1737' }
    @{ Cell = 'E10'; Value = '3.51%' }
    @{ Cell = 'D11'; Value = '0.07255' }
    @{ Cell = 'E11'; Value = '2.00%' }
    @{ Cell = 'D12'; Value = '0.08331' }
    @{ Cell = 'E12'; Value = '3.65%' }
    @{ Cell = 'D13'; Value = '0.02995' }
    @{ Cell = 'E13'; Value = '0.41%' }
    @{ Cell = 'D14'; Value = '0.09954' }
    @{ Cell = 'E14'; Value = '0.54%' }
    @{ Cell = 'D15'; Value = '0.001493' }
    @{ Cell = 'D16'; Value = '0.006109' }
    @{ Cell = 'E16'; Value = '-1.17%' }
    @{ Cell = 'D18'; Value = '2.252' }
    @{ Cell = 'E18'; Value = '0.92%' }
    @{ Cell = 'E19'; Value = '0.35%' }
    @{ Cell = 'E20'; Value = '-0.35%' }
    @{ Cell = 'D21'; Value = '4.635' }
    @{ Cell = 'E21'; Value = '1.81%' }
    @{ Cell = 'E22'; Value = '3.33%' }
    @{ Cell = 'D23'; Value = '0.04583' }
    @{ Cell = 'E23'; Value = '-1.17%' }
    @{ Cell = 'D24'; Value = '0.001259' }
    @{ Cell = 'E24'; Value = '3.51%' }
    @{ Cell = 'D25'; Value = '0.004444' }
    @{ Cell = 'E25'; Value = '0.35%' }
    @{ Cell = 'D26'; Value = '0.0001181' }
    @{ Cell = 'E26'; Value = '-8.98%' }
    @{ Cell = 'D27'; Value = '0.0003433' }
    @{ Cell = 'E27'; Value = '83.40%' }
    @{ Cell = 'D39'; Value = '0.01825' }
    @{ Cell = 'E39'; Value = '7.40%' }
    @{ Cell = 'D40'; Value = '0.04513' }
    @{ Cell = 'E40'; Value = '2.22%' }
    @{ Cell = 'D41'; Value = '0.007019' }
    @{ Cell = 'E41'; Value = '-2.66%' }
    @{ Cell = 'D42'; Value = '0.1343' }
    @{ Cell = 'E42'; Value = '1.04%' }
    @{ Cell = 'D43'; Value = '0.002242' }
    @{ Cell = 'E43'; Value = '4.96%' }
    @{ Cell = 'D44'; Value = '0.009832' }
    @{ Cell = 'E44'; Value = '-10.57%' }
    @{ Cell = 'D45'; Value = '0.00006475' }
    @{ Cell = 'E45'; Value = '7.94%' }
    @{ Cell = 'D46'; Value = '0.00000000751' }
    @{ Cell = 'E46'; Value = '0.09%' }
    @{ Cell = 'D47'; Value = '0.006206' }
    @{ Cell = 'E47'; Value = '-39.23%' }
    @{ Cell = 'E48'; Value = '-56.68%' }
    @{ Cell = 'D49'; Value = '0.00002102' }
    @{ Cell = 'E49'; Value = '0.09%' }
    @{ Cell = 'D50'; Value = '0.0002002' }
    @{ Cell = 'E50'; Value = '0.16%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
